$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1453.3572
$ws.Range("I19").Value = 1425.5
$ws.Range("J19").Value = 1523
$ws.Range("K19").Value = 1425.5
$ws.Range("L19").Value = 1523
$ws.Range("M19").Value = -1250.5
$ws.Range("N19").Value = -1873
$ws.Range("H76").Value = 4439.2
$ws.Range("I76").Value = 4303.6
$ws.Range("J76").Value = 4574.8
$ws.Range("K76").Value = 4303.6
$ws.Range("L76").Value = 4574.8
$ws.Range("M76").Value = -3988.6
$ws.Range("N76").Value = -5204.8
$ws.Range("H79").Value = 4439.2
$ws.Range("I79").Value = 4303.6
$ws.Range("J79").Value = 4574.8
$ws.Range("K79").Value = 4303.6
$ws.Range("L79").Value = 4574.8
$ws.Range("M79").Value = -3211.6
$ws.Range("N79").Value = -6758.8
$ws.Range("H106").Value = 5620
$ws.Range("I106").Value = 3799.5
$ws.Range("J106").Value = 7440.5
$ws.Range("K106").Value = 3799.5
$ws.Range("L106").Value = 7440.5
$ws.Range("M106").Value = -3168.5
$ws.Range("N106").Value = -8702.5
$ws.Range("H116").Value = 3772.1035
$ws.Range("I116").Value = 3858.1052
$ws.Range("J116").Value = 3608.7
$ws.Range("K116").Value = 3858.1052
$ws.Range("L116").Value = 3608.7
$ws.Range("M116").Value = -416.1052
$ws.Range("N116").Value = -10492.7
$ws.Range("H135").Value = 50000732
$ws.Range("I135").Value = 50000732
$ws.Range("K135").Value = 450006588
$ws.Range("M135").Value = -450004053
$ws.Range("H137").Value = 2580.8076
$ws.Range("I137").Value = 2184
$ws.Range("J137").Value = 3473.625
$ws.Range("K137").Value = 6552
$ws.Range("L137").Value = 10420.875
$ws.Range("M137").Value = -4002
$ws.Range("N137").Value = -15520.875
$ws.Range("H138").Value = 2540.9836
$ws.Range("I138").Value = 1393.7931
$ws.Range("K138").Value = 4181.379300000001
$ws.Range("M138").Value = 958.6206999999995
$ws.Range("H141").Value = 930.1724
$ws.Range("I141").Value = 965
$ws.Range("K141").Value = 2895
$ws.Range("M141").Value = 2285

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3468.3428
$ws.Range("I32").Value = 3296.4
$ws.Range("J32").Value = 3898.2
$ws.Range("K32").Value = 3296.4
$ws.Range("L32").Value = 3898.2
$ws.Range("M32").Value = -3009.4
$ws.Range("N32").Value = -4472.2
$ws.Range("H61").Value = 66668000
$ws.Range("I61").Value = 76924110
$ws.Range("K61").Value = 76924110
$ws.Range("M61").Value = -76923898
$ws.Range("H74").Value = 37042150
$ws.Range("I74").Value = 41671800
$ws.Range("J74").Value = 5000
$ws.Range("K74").Value = 41671800
$ws.Range("L74").Value = 5000
$ws.Range("M74").Value = -41670926
$ws.Range("N74").Value = -6748
$ws.Range("H77").Value = 37042150
$ws.Range("I77").Value = 41671800
$ws.Range("J77").Value = 5000
$ws.Range("K77").Value = 208359000
$ws.Range("L77").Value = 25000
$ws.Range("M77").Value = -208354632
$ws.Range("N77").Value = -33736
$ws.Range("H132").Value = 3128512.5
$ws.Range("I132").Value = 3849733
$ws.Range("J132").Value = 3223
$ws.Range("K132").Value = 11549199
$ws.Range("L132").Value = 9669
$ws.Range("M132").Value = -11546669
$ws.Range("N132").Value = -14729
$ws.Range("H136").Value = 66668000
$ws.Range("I136").Value = 76924110
$ws.Range("K136").Value = 230772330
$ws.Range("M136").Value = -230769780

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 26842864
$ws.Range("I134").Value = 26842864
$ws.Range("K134").Value = 80528592
$ws.Range("M134").Value = -80526057

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7439.025
$ws.Range("I31").Value = 4550.28
$ws.Range("J31").Value = 12253.6
$ws.Range("K31").Value = 4550.28
$ws.Range("L31").Value = 12253.6
$ws.Range("M31").Value = -4255.28
$ws.Range("N31").Value = -12843.6
$ws.Range("H34").Value = 7439.025
$ws.Range("I34").Value = 4550.28
$ws.Range("J34").Value = 12253.6
$ws.Range("K34").Value = 4550.28
$ws.Range("L34").Value = 12253.6
$ws.Range("M34").Value = -4348.28
$ws.Range("N34").Value = -12657.6
$ws.Range("H58").Value = 20005222
$ws.Range("I58").Value = 27784314
$ws.Range("J58").Value = 1844.1428
$ws.Range("K58").Value = 27784314
$ws.Range("L58").Value = 1844.1428
$ws.Range("M58").Value = -27784111
$ws.Range("N58").Value = -2250.1428
$ws.Range("H132").Value = 19609606
$ws.Range("I132").Value = 20835136
$ws.Range("J132").Value = 1116.3334
$ws.Range("K132").Value = 62505408
$ws.Range("L132").Value = 3349.0002
$ws.Range("M132").Value = -62502878
$ws.Range("N132").Value = -8409.0002
$ws.Range("H134").Value = 6251960.5
$ws.Range("I134").Value = 7577674
$ws.Range("J134").Value = 2170
$ws.Range("K134").Value = 22733022
$ws.Range("L134").Value = 6510
$ws.Range("M134").Value = -22730487
$ws.Range("N134").Value = -11580
$ws.Range("H136").Value = 20005222
$ws.Range("I136").Value = 27784314
$ws.Range("J136").Value = 1844.1428
$ws.Range("K136").Value = 83352942
$ws.Range("L136").Value = 5532.428400000001
$ws.Range("M136").Value = -83350392
$ws.Range("N136").Value = -10632.4284
$ws.Range("H141").Value = 90897.5
$ws.Range("J141").Value = 90897.5
$ws.Range("L141").Value = 90897.5
$ws.Range("N141").Value = -101257.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 143793.28
$ws.Range("I113").Value = 200770.8
$ws.Range("J113").Value = 1349.5
$ws.Range("K113").Value = 602312.3999999999
$ws.Range("L113").Value = 4048.5
$ws.Range("M113").Value = -600142.3999999999
$ws.Range("N113").Value = -8388.5
$ws.Range("H122").Value = 549.0476
$ws.Range("J122").Value = 977.1429000000001
$ws.Range("L122").Value = 8794.286100000001
$ws.Range("N122").Value = -13694.2861

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 5032.7812
$ws.Range("I122").Value = 3331
$ws.Range("J122").Value = 8776.700000000001
$ws.Range("K122").Value = 9993
$ws.Range("L122").Value = 26330.1
$ws.Range("M122").Value = -7543
$ws.Range("N122").Value = -31230.1

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4153.4614
$ws.Range("I7").Value = 4141
$ws.Range("K7").Value = 4141
$ws.Range("M7").Value = -4029
$ws.Range("H46").Value = 979.5454999999999
$ws.Range("I46").Value = 999
$ws.Range("J46").Value = 972.25
$ws.Range("K46").Value = 999
$ws.Range("L46").Value = 972.25
$ws.Range("M46").Value = -811
$ws.Range("N46").Value = -1348.25
$ws.Range("H100").Value = 19448660
$ws.Range("I100").Value = 29170490
$ws.Range("J100").Value = 5000
$ws.Range("K100").Value = 29170490
$ws.Range("L100").Value = 5000
$ws.Range("M100").Value = -29169949
$ws.Range("N100").Value = -6082
$ws.Range("H104").Value = 21671.25
$ws.Range("J104").Value = 21671.25
$ws.Range("L104").Value = 21671.25
$ws.Range("N104").Value = -28659.25
$ws.Range("H126").Value = 4153.4614
$ws.Range("I126").Value = 4141
$ws.Range("K126").Value = 12423
$ws.Range("M126").Value = -9953
$ws.Range("H132").Value = 15781769
$ws.Range("I132").Value = 16908758
$ws.Range("J132").Value = 3933.3333
$ws.Range("K132").Value = 50726274
$ws.Range("L132").Value = 11799.9999
$ws.Range("M132").Value = -50723744
$ws.Range("N132").Value = -16859.9999
$ws.Range("H136").Value = 1887.2941
$ws.Range("I136").Value = 955.375
$ws.Range("J136").Value = 2715.6667
$ws.Range("K136").Value = 2866.125
$ws.Range("L136").Value = 8147.000100000001
$ws.Range("M136").Value = -316.125
$ws.Range("N136").Value = -13247.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1245
$ws.Range("I113").Value = 1337.3572
$ws.Range("J113").Value = 814
$ws.Range("K113").Value = 4012.0716
$ws.Range("L113").Value = 2442
$ws.Range("M113").Value = -1842.0716
$ws.Range("N113").Value = -6782
$ws.Range("H126").Value = 1192.95
$ws.Range("I126").Value = 1021.1177
$ws.Range("K126").Value = 3063.3531
$ws.Range("M126").Value = -593.3531000000003
$ws.Range("H132").Value = 9260516
$ws.Range("I132").Value = 10205321
$ws.Range("K132").Value = 30615963
$ws.Range("M132").Value = -30613433
